$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (G1, H1) - same style as existing header row (s="1")
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Updated values in existing columns (B, C, D) for rows 2 and 3
$ws.Range("B2").Value = 0.06827263118387658
$ws.Range("C2").Value = 0.9987528207119
$ws.Range("D2").Value = 0.206831987943706

$ws.Range("B3").Value = 0.1253955732895198
$ws.Range("C3").Value = 0.9907353182326988
$ws.Range("D3").Value = 0.276136362844067

# New values in columns G and H
$ws.Range("G2").Value = 0.2668650318499809
$ws.Range("H2").Value = 0.998

$ws.Range("G3").Value = 0.2668650318499809
$ws.Range("H3").Value = 0.998
